# "Agregada función de dias para fecha límite"
#
# In the "Entidad Tarea" table (Tabla4, A22:E28) on sheet "Hoja1", the
# attribute row that used to describe the "Descripción" of the task's
# function is repurposed into a new "Fase" attribute row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 26: change "Descripción" / "Descripción de la función" / "Terminar el
# archivo XXXX..." into the new "Fase" attribute.
$ws.Range("A26").Value = "Fase"
$ws.Range("B26").Value = "Fase en la que se ecuentra la tarea"
$ws.Range("C26").Value = "Caracteres(100)"
$ws.Range("E26").Value = "El archivo XXXX se encuentra en Fase YYYY"

# The new, shorter text needs fewer wrapped lines than the old one.
$ws.Rows(26).RowHeight = 41.4

# Leave the workbook scrolled to / selecting the edited cell, as in the
# authored change.
$ws.Range("E26").Select()
